# Table_01_AdonisRes_Taxonomic_19-04-20.docx - taxonomic revision
#
# The table values were revised with extra decimal precision. Each touched
# cell originally held a single run of text (e.g. "0.212"); after the edit
# the cell holds the same visible number but built from TWO runs (e.g.
# "0.2" + "04" = "0.204"), mirroring how the authors' tool re-typed the
# trailing digits of each figure. A couple of cells (bold, statistically
# significant p-values) keep their bold run-properties on both runs.

function Split-TableCell {
    param($table, $row, $col, $firstText, $secondText, $isBold)

    $cellRange = $table.Cell($row, $col).Range
    # Trim the end-of-cell mark (wdCharacter unit collapses the 2 hidden
    # end-of-cell chars into a single logical "character" move).
    $cellRange.MoveEnd(1, -1)

    # Replace the cell's whole text with the first run's text.
    $run1 = $word.ActiveDocument.Range($cellRange.Start, $cellRange.End)
    $run1.Text = $firstText
    # Touch Bold (toggle then set) so the engine treats this as a distinct
    # run even when its final formatting matches its neighbour.
    $run1.Bold = 1
    if ($isBold -eq 1) {
        $run1.Bold = 1
    } else {
        $run1.Bold = 0
    }

    # Insert the second run immediately after the first.
    $run2 = $word.ActiveDocument.Range($run1.End, $run1.End)
    $run2.InsertAfter($secondText)
    $run2.Bold = 1
    if ($isBold -eq 1) {
        $run2.Bold = 1
    } else {
        $run2.Bold = 0
    }
}

$d = $word.ActiveDocument
$t = $d.Tables.Item(1)

# Taxonomic / Total / R2 (Turbidity row)
Split-TableCell $t 4 4 "0.2" "04" 0
# Taxonomic / Turnover / R2 (Turbidity row)
Split-TableCell $t 4 6 "0.06" "8" 0
# Taxonomic / Turnover / P value (Turbidity row)
Split-TableCell $t 4 7 "0." "094" 0
# Taxonomic / Nestedness / R2 (Turbidity row)
Split-TableCell $t 4 9 "-0.06" "3" 0
# Taxonomic / Nestedness / P value (Turbidity row)
Split-TableCell $t 4 10 "0.82" "9" 0

# Taxonomic / Total / P value (pH row) - bold (significant)
Split-TableCell $t 5 4 "0.0" "47" 1
# Taxonomic / Turnover / P value (pH row)
Split-TableCell $t 5 7 "0.0" "67" 0
# Taxonomic / Nestedness / R2 (pH row)
Split-TableCell $t 5 9 "0.34" "9" 0
# Taxonomic / Nestedness / P value (pH row)
Split-TableCell $t 5 10 "0.09" "2" 0

# Taxonomic / Total / P value (Conductivity row)
Split-TableCell $t 6 4 "0." "299" 0
# Taxonomic / Turnover / P value (Conductivity row)
Split-TableCell $t 6 7 "0." "197" 0
# Taxonomic / Nestedness / R2 (Conductivity row)
Split-TableCell $t 6 9 "-0.03" "6" 0

# Taxonomic / Total / P value (Dissolved O2 row)
Split-TableCell $t 7 4 "0.8" "14" 0
# Taxonomic / Turnover / P value (Dissolved O2 row)
Split-TableCell $t 7 7 "0.5" "68" 0
# Taxonomic / Nestedness / P value (Dissolved O2 row)
Split-TableCell $t 7 10 "0.6" "41" 0

# Taxonomic / Total / P value (Altitude row)
Split-TableCell $t 8 4 "0.6" "25" 0
# Taxonomic / Turnover / R2 (Altitude row)
Split-TableCell $t 8 6 "-0.00" "2" 0
# Taxonomic / Turnover / P value (Altitude row)
Split-TableCell $t 8 7 "0.9" "19" 0
# Taxonomic / Nestedness / R2 (Altitude row)
Split-TableCell $t 8 9 "0.0" "89" 0
# Taxonomic / Nestedness / P value (Altitude row)
Split-TableCell $t 8 10 "0." "323" 0

# Taxonomic / Turnover / P value (Temperature row)
Split-TableCell $t 9 7 "0.2" "72" 0
# Taxonomic / Nestedness / P value (Temperature row)
Split-TableCell $t 9 10 "0.5" "12" 0

# Taxonomic / Total / R2 (Depth row)
Split-TableCell $t 10 4 "0.6" "20" 0
# Taxonomic / Turnover / P value (Depth row)
Split-TableCell $t 10 7 "0.7" "31" 0
# Taxonomic / Nestedness / R2 (Depth row)
Split-TableCell $t 10 9 "0.03" "4" 0
# Taxonomic / Nestedness / P value (Depth row)
Split-TableCell $t 10 10 "0.4" "61" 0

# Taxonomic / Total / P value (Water velocity row)
Split-TableCell $t 11 4 "0.09" "2" 0
# Taxonomic / Turnover / R2 (Water velocity row)
Split-TableCell $t 11 6 "0.08" "8" 0
# Taxonomic / Turnover / P value (Water velocity row)
Split-TableCell $t 11 7 "0.0" "56" 0
# Taxonomic / Nestedness / R2 (Water velocity row) - simple text shrink, no split
$cell11_9 = $t.Cell(11, 9).Range
$cell11_9.MoveEnd(1, -1)
$cell11_9.Text = "-0.001"
# Taxonomic / Nestedness / P value (Water velocity row)
Split-TableCell $t 11 10 "0.6" "41" 0

# Taxonomic / Total / P value (Basin row)
Split-TableCell $t 12 4 "0.0" "73" 0
# Taxonomic / Turnover / P value (Basin row) - bold (significant)
Split-TableCell $t 12 7 "0.0" "39" 1
# Taxonomic / Nestedness / R2 (Basin row)
Split-TableCell $t 12 9 "-0.13" "1" 0
# Taxonomic / Nestedness / P value (Basin row)
Split-TableCell $t 12 10 "0.8" "79" 0

# Taxonomic / Total / R2 (Residual row)
Split-TableCell $t 13 3 "0.35" "1" 0
# Taxonomic / Turnover / R2 (Residual row)
Split-TableCell $t 13 6 "0.31" "1" 0
